# Insert two new daily-price rows for "Vega Monumental Concepción" (Limón)
# right above the existing row 797, shifting the rest of the table down by
# two rows (old row 797 -> new row 799, ... old row 891 -> new row 893).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 797.. down by two rows, leaving 797:798 empty for the new data.
$ws.Rows("797:798").Insert()

# --- New row 797 -----------------------------------------------------
$ws.Range("A797").Value = 11
$ws.Range("B797").Value = "Vega Monumental Concepción"
$ws.Range("C797").Value = "Bíobío"
$ws.Range("D797").Value = 45212
$ws.Range("D797").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E797").Value = 8
$ws.Range("F797").Value = "Fruta"
$ws.Range("G797").Value = 100102
$ws.Range("H797").Value = "Cítricos"
$ws.Range("I797").Value = 100102003
$ws.Range("J797").Value = "Limón"
$ws.Range("K797").Value = "Sin especificar"
$ws.Range("L797").Value = "1a amarillo"
$ws.Range("M797").Value = 300
$ws.Range("N797").Value = 9000
$ws.Range("O797").Value = 10000
$ws.Range("P797").Value = 9667
$ws.Range("Q797").Value = '$/malla 18 kilos'
$ws.Range("R797").Value = "Provincia de Melipilla"
$ws.Range("S797").Value = 537
$ws.Range("T797").Value = 18

# --- New row 798 -----------------------------------------------------
$ws.Range("A798").Value = 11
$ws.Range("B798").Value = "Vega Monumental Concepción"
$ws.Range("C798").Value = "Bíobío"
$ws.Range("D798").Value = 45212
$ws.Range("D798").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E798").Value = 8
$ws.Range("F798").Value = "Fruta"
$ws.Range("G798").Value = 100102
$ws.Range("H798").Value = "Cítricos"
$ws.Range("I798").Value = 100102003
$ws.Range("J798").Value = "Limón"
$ws.Range("K798").Value = "Sin especificar"
$ws.Range("L798").Value = "2a amarillo"
$ws.Range("M798").Value = 100
$ws.Range("N798").Value = 8000
$ws.Range("O798").Value = 8000
$ws.Range("P798").Value = 8000
$ws.Range("Q798").Value = '$/malla 18 kilos'
$ws.Range("R798").Value = "Provincia de Melipilla"
$ws.Range("S798").Value = 444
$ws.Range("T798").Value = 18
